# Updates the cryptos list data (prices and 1h volume deltas) to reflect
# the latest scrape, per the "Updated cryptos list ... with GitHub
# Actions" commit.
#
# Numeric-looking price strings (e.g. "1.00", "14.30") must stay as TEXT
# (the Price column stores values as inline strings, not numbers), so a
# leading apostrophe is used to stop Excel auto-converting them to
# numbers (which would drop things like trailing zeros). The cell's
# Style is then reset back to "Normal" so no lingering quote-prefix /
# text-format style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.004.92"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "2.023.20"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "'225.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("D6").Value = "'0.606"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.84%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'54.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.38%  "
$ws.Range("D9").Value = "'0.379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").Value = "'0.0787"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("D12").Value = "2.315.44"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("D13").Value = "'14.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("D14").Value = "'20.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "'0.746"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "'5.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.29%  "
$ws.Range("D17").Value = "2.020.30"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "36.921.15"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "'6.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.73%  "
$ws.Range("D20").Value = "'68.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "0.0₃0821"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "'226.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("E25").Value = "  -7.31%  "
$ws.Range("D26").Value = "'165.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "'9.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.10%  "
$ws.Range("D28").Value = "'18.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.04%  "
$ws.Range("D29").Value = "'0.124"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.06%  "
$ws.Range("D30").Value = "'1.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("D31").Value = "'0.117"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.88%  "
$ws.Range("D32").Value = "'4.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("D33").Value = "'0.0617"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").Value = "'4.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("D35").Value = "'2.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.59%  "
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "'0.996"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'3.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.81%  "
$ws.Range("D39").Value = "'5.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'17.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.19%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.490.41"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "'0.0218"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.17%  "
$ws.Range("D43").Value = "'95.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("D44").Value = "'0.0926"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("E45").Value = "  -5.41%  "
$ws.Range("E46").Value = "  -5.40%  "
$ws.Range("D47").Value = "'7.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.55%  "
$ws.Range("D49").Value = "'2.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "2.207.24"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("D51").Value = "'3.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.11%  "
